$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.811.54"
$ws.Range("E2").Value = "  +1.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.514.44"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.77"
$ws.Range("E5").Value = "  +3.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.11"
$ws.Range("E6").Value = "  -2.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.615"
$ws.Range("E7").Value = "  -1.59%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.510.93"
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").Value = "  +4.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.69"
$ws.Range("E11").Value = "  -1.17%  "
$ws.Range("E12").Value = "  -3.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "47.20"
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("E14").Value = "  +0.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.084.75"
$ws.Range("E15").Value = "  +0.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "619.96"
$ws.Range("E16").Value = "  -8.53%  "
$ws.Range("E17").Value = "  -4.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.510.55"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.820.15"
$ws.Range("E19").Value = "  +1.04%  "
$ws.Range("E20").Value = "  -2.10%  "
$ws.Range("E21").Value = "  -1.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.98"
$ws.Range("E22").Value = "  -11.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.72"
$ws.Range("E24").Value = "  -2.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "95.82"
$ws.Range("E25").Value = "  -2.56%  "
$ws.Range("E26").Value = "  -0.80%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("E28").Value = "  -2.97%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.21"
$ws.Range("E29").Value = "  -3.09%  "
$ws.Range("E30").Value = "  +0.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.42"
$ws.Range("E31").Value = "  -4.19%  "
$ws.Range("E32").Value = "  -4.82%  "
$ws.Range("E33").Value = "  -1.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.96"
$ws.Range("E34").Value = "  -5.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "568.09"
$ws.Range("E35").Value = "  -2.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.76"
$ws.Range("E36").Value = "  -1.62%  "
$ws.Range("E37").Value = "  -3.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "57.04"
$ws.Range("E38").Value = "  -0.59%  "
$ws.Range("E39").Value = "  -4.07%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("E41").Value = "  +2.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0448"
$ws.Range("E42").Value = "  +1.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.326"
$ws.Range("E43").Value = "  -4.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.333.35"
$ws.Range("E44").Value = "  -2.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.00"
$ws.Range("E45").Value = "  +2.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "33.09"
$ws.Range("E46").Value = "  -1.40%  "
$ws.Range("E47").Value = "  -1.21%  "
$ws.Range("E48").Value = "  +0.87%  "
$ws.Range("E49").Value = "  -3.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "135.66"
$ws.Range("E50").Value = "  +2.64%  "
$ws.Range("E51").Value = "  +0.06%  "
